# "Lista de chequeo version dos"
# - Insert a new (blank) column before column C, pushing the existing
#   C:H content to D:I.
# - Resize the (now) D:I columns to their new widths.
# - Update the text of the "ID 07" question.
# - Let row heights follow the new wrapping (row 4 and rows 8-18 change).
# - Update the active selection / view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a column before C (shifts C:H -> D:I) -------------------------
$ws.Columns.Item(3).Insert()

# --- New column widths (D:I) ----------------------------------------------
# NOTE: Excel's ColumnWidth (in characters) is stored in the worksheet XML
# as ColumnWidth + 5/6, so we back the target stored widths out by 5/6.
$ws.Columns.Item(4).ColumnWidth = 13.7109375 - 0.8333333333333334
$ws.Columns.Item(5).ColumnWidth = 30.5703125 - 0.8333333333333334
$ws.Columns.Item(6).ColumnWidth = 14        - 0.8333333333333334
$ws.Columns.Item(7).ColumnWidth = 15.140625 - 0.8333333333333334
$ws.Columns.Item(8).ColumnWidth = 23.7109375 - 0.8333333333333334
$ws.Columns.Item(9).ColumnWidth = 24.28515625 - 0.8333333333333334

# --- Update the "ID 07" question text (now in E14) -------------------------
$ws.Range("E14").Value = "¿Como cajero tiene un buena organizacion en su espacio?"

# --- Row heights (reflow after the new column widths/word-wrap) -----------
$ws.Rows.Item(4).RowHeight = 39.75
$ws.Rows.Item(5).RowHeight = 31.5
$ws.Rows.Item(8).RowHeight = 38.25
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 42.75
$ws.Rows.Item(11).RowHeight = 42
$ws.Rows.Item(12).RowHeight = 33.75
$ws.Rows.Item(13).RowHeight = 40.5
$ws.Rows.Item(14).RowHeight = 33.75
$ws.Rows.Item(15).RowHeight = 46.5
$ws.Rows.Item(16).RowHeight = 35.25
$ws.Rows.Item(17).RowHeight = 43.5
$ws.Rows.Item(18).RowHeight = 54.75

# --- View / selection -------------------------------------------------------
$ws.Range("O10").Select()
